$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.824.95"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.642.61"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.07"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.27"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("E10").Value = "  +5.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.59"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.48"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.114.30"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.660.28"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.656.23"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.80"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.57"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.62"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.92"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.31"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.67"
$ws.Range("E25").Value = "  +7.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.70"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  +5.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "567.77"
$ws.Range("E28").Value = "  +4.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.10"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0850"
$ws.Range("E33").Value = "  +4.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.95"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.405"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.96"
$ws.Range("E39").Value = "  +4.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.19"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "165.21"
$ws.Range("E42").Value = "  -6.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.17"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.79"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.94"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0567"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.626"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("E48").Value = "  +14.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0246"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.85"
$ws.Range("E51").Value = "  -0.68%  "
